$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D6").Value = "[Object Detection] YOLO v1 ~ v6 비교(1)"
$ws.Range("E6").Value = "https://leedakyeong.tistory.com/entry/Object-Detection-YOLO-v1v6-%EB%B9%84%EA%B5%90"

$ws.Range("D9").Value = "전공자 vs. 비전공자 in Data Science"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/majored-vs-not-in-data-science/#utm_source=rss&utm_medium=rss&utm_campaign=majored-vs-not-in-data-science"

$ws.Range("D51").Value = "[vim] vim의 복사, 붙여넣기, 잘라내기 단축키"
$ws.Range("E51").Value = "https://bskyvision.com/1309"
